$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.998.09'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '2.966.66'
$ws.Range("E3").Value = '  +2.98%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''201.49'
$ws.Range("E5").Value = '  +2.71%  '
$ws.Range("D6").Value = '''598.41'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.550'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("D10").Value = '2.963.96'
$ws.Range("E10").Value = '  +2.88%  '
$ws.Range("E11").Value = '  +13.80%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").Value = '3.508.59'
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("D14").Value = '''4.91'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '76.873.04'
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("D16").Value = '''28.38'
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").Value = '2.953.83'
$ws.Range("E18").Value = '  +2.90%  '
$ws.Range("D19").Value = '''13.53'
$ws.Range("E19").Value = '  +7.78%  '
$ws.Range("D20").Value = '''8.72'
$ws.Range("D21").Value = '''373.99'
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").Value = '''4.33'
$ws.Range("E22").Value = '  +4.65%  '
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").Value = '''72.78'
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("D25").Value = '3.114.87'
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").Value = '''4.28'
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("D28").Value = '''9.73'
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("D29").Value = '''0.0000108'
$ws.Range("E29").Value = '  +3.14%  '
$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '''8.30'
$ws.Range("E31").Value = '  +6.88%  '
$ws.Range("D33").Value = '''499.36'
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("D34").Value = '''1.84'
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '''166.31'
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '''0.399'
$ws.Range("E37").Value = '  +15.68%  '
$ws.Range("E38").Value = '  +22.48%  '
$ws.Range("D39").Value = '''20.22'
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("D40").Value = '''19.80'
$ws.Range("E40").Value = '  +1.36%  '
$ws.Range("D41").Value = '''0.111'
$ws.Range("E41").Value = '  -6.17%  '
$ws.Range("D43").Value = '''181.35'
$ws.Range("E43").Value = '  -2.13%  '
$ws.Range("D44").Value = '''4.94'
$ws.Range("E44").Value = '  -2.68%  '
$ws.Range("D45").Value = '''1.65'
$ws.Range("E45").Value = '  -1.54%  '
$ws.Range("D46").Value = '''40.12'
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").Value = '''1.19'
$ws.Range("E47").Value = '  -3.42%  '
$ws.Range("D48").Value = '''0.593'
$ws.Range("E48").Value = '  +2.22%  '
$ws.Range("D49").Value = '''3.90'
$ws.Range("E49").Value = '  +3.92%  '
$ws.Range("D50").Value = '''2.33'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("D51").Value = '''22.74'
$ws.Range("E51").Value = '  +5.14%  '
